# Updated cryptos list on Mon Jul 17 05:35:47 UTC 2023 with GitHub Actions
#
# Refreshes the scraped Price / Volume(1h) columns for every coin row with
# the latest values. Two pairs of rows (46/47 and 49/50) also changed which
# coin they describe (EnergySwap <-> Aptos, Elrond <-> Algorand swapped
# positions in the source ranking), so their Coin name and Link columns are
# updated as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new text value, exactly as scraped.
$updates = @(
    @{ Cell = "D2"; Value = '30.324.12' }
    @{ Cell = "E2"; Value = '  +0.28%  ' }
    @{ Cell = "D3"; Value = '1.934.79' }
    @{ Cell = "E3"; Value = '  +0.42%  ' }
    @{ Cell = "E4"; Value = '  -0.06%  ' }
    @{ Cell = "D5"; Value = '0.7502' }
    @{ Cell = "E5"; Value = '  +5.34%  ' }
    @{ Cell = "D6"; Value = '243.33' }
    @{ Cell = "E6"; Value = '  -2.13%  ' }
    @{ Cell = "D7"; Value = '1.001' }
    @{ Cell = "D8"; Value = '28.13' }
    @{ Cell = "E8"; Value = '  +2.87%  ' }
    @{ Cell = "D9"; Value = '0.3188' }
    @{ Cell = "E9"; Value = '  -0.44%  ' }
    @{ Cell = "D10"; Value = '0.07042' }
    @{ Cell = "E10"; Value = '  -0.24%  ' }
    @{ Cell = "D11"; Value = '0.7831' }
    @{ Cell = "E11"; Value = '  -1.05%  ' }
    @{ Cell = "D12"; Value = '0.08047' }
    @{ Cell = "E12"; Value = '  +1.12%  ' }
    @{ Cell = "D13"; Value = '1.938.69' }
    @{ Cell = "E13"; Value = '  +0.51%  ' }
    @{ Cell = "D14"; Value = '5.405' }
    @{ Cell = "E14"; Value = '  +0.50%  ' }
    @{ Cell = "D15"; Value = '93.13' }
    @{ Cell = "E15"; Value = '  -1.79%  ' }
    @{ Cell = "D16"; Value = '14.53' }
    @{ Cell = "E16"; Value = '  -0.51%  ' }
    @{ Cell = "D17"; Value = '30.331.45' }
    @{ Cell = "E17"; Value = '  +0.20%  ' }
    @{ Cell = "D18"; Value = '6.108' }
    @{ Cell = "E18"; Value = '  +6.15%  ' }
    @{ Cell = "D19"; Value = '252.65' }
    @{ Cell = "E19"; Value = '  -2.04%  ' }
    @{ Cell = "D20"; Value = '0.000008005' }
    @{ Cell = "E20"; Value = '  -0.26%  ' }
    @{ Cell = "D21"; Value = '2.182.09' }
    @{ Cell = "E21"; Value = '  -0.03%  ' }
    @{ Cell = "D22"; Value = '0.9993' }
    @{ Cell = "E22"; Value = '  -0.14%  ' }
    @{ Cell = "E23"; Value = '  -0.09%  ' }
    @{ Cell = "D24"; Value = '6.713' }
    @{ Cell = "E24"; Value = '  -1.89%  ' }
    @{ Cell = "D25"; Value = '9.582' }
    @{ Cell = "E25"; Value = '  +0.58%  ' }
    @{ Cell = "D26"; Value = '164.99' }
    @{ Cell = "E26"; Value = '  -0.59%  ' }
    @{ Cell = "D27"; Value = '19.09' }
    @{ Cell = "E27"; Value = '  -0.05%  ' }
    @{ Cell = "D28"; Value = '0.1309' }
    @{ Cell = "E28"; Value = '  +4.08%  ' }
    @{ Cell = "D29"; Value = '2.207' }
    @{ Cell = "E29"; Value = '  -2.42%  ' }
    @{ Cell = "D30"; Value = '1.368' }
    @{ Cell = "E30"; Value = '  +0.95%  ' }
    @{ Cell = "D31"; Value = '1.542' }
    @{ Cell = "E31"; Value = '  +0.94%  ' }
    @{ Cell = "D32"; Value = '4.441' }
    @{ Cell = "E32"; Value = '  +1.19%  ' }
    @{ Cell = "D33"; Value = '4.153' }
    @{ Cell = "E33"; Value = '  +0.76%  ' }
    @{ Cell = "D34"; Value = '0.05289' }
    @{ Cell = "E34"; Value = '  +2.89%  ' }
    @{ Cell = "D35"; Value = '1.337' }
    @{ Cell = "E35"; Value = '  +5.35%  ' }
    @{ Cell = "D36"; Value = '0.7580' }
    @{ Cell = "E36"; Value = '  +1.67%  ' }
    @{ Cell = "D37"; Value = '2.781' }
    @{ Cell = "E37"; Value = '  +0.71%  ' }
    @{ Cell = "D38"; Value = '0.01966' }
    @{ Cell = "E38"; Value = '  +0.28%  ' }
    @{ Cell = "D39"; Value = '2.801' }
    @{ Cell = "E39"; Value = '  +0.13%  ' }
    @{ Cell = "D40"; Value = '79.09' }
    @{ Cell = "E40"; Value = '  +2.14%  ' }
    @{ Cell = "D41"; Value = '6.536' }
    @{ Cell = "E41"; Value = '  +2.89%  ' }
    @{ Cell = "D42"; Value = '0.4515' }
    @{ Cell = "E42"; Value = '  +0.52%  ' }
    @{ Cell = "D43"; Value = '1.986' }
    @{ Cell = "E43"; Value = '  +0.12%  ' }
    @{ Cell = "D44"; Value = '1.001' }
    @{ Cell = "D45"; Value = '0.8401' }
    @{ Cell = "E45"; Value = '  -0.61%  ' }
    @{ Cell = "B46"; Value = 'Aptos' }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = "D46"; Value = '7.722' }
    @{ Cell = "E46"; Value = '  +4.04%  ' }
    @{ Cell = "B47"; Value = 'EnergySwap' }
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = "D47"; Value = '10.03' }
    @{ Cell = "E47"; Value = '  +3.18%  ' }
    @{ Cell = "D48"; Value = '101.87' }
    @{ Cell = "E48"; Value = '  +1.36%  ' }
    @{ Cell = "B49"; Value = 'Algorand' }
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' }
    @{ Cell = "D49"; Value = '0.1263' }
    @{ Cell = "E49"; Value = '  +11.28%  ' }
    @{ Cell = "B50"; Value = 'Elrond' }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' }
    @{ Cell = "D50"; Value = '37.69' }
    @{ Cell = "E50"; Value = '  +3.16%  ' }
    @{ Cell = "D51"; Value = '968.51' }
    @{ Cell = "E51"; Value = '  +5.77%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force Text format BEFORE assigning so Excel keeps the original string
    # representation (e.g. "30.324.12", "0.7502", "  +0.28%  ") instead of
    # silently reinterpreting numeric-looking text as a real number.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Reset to the default style so no stray per-cell number format lingers,
    # matching the original plain (unstyled) data cells.
    $cell.Style = "Normal"
}
